# Stable release 1.8 FINAL_FINAL_FINAL
#
# Remove the duplicate "Stake holders" title slide (originally slide position 3)
# and the trailing empty slide (originally slide position 12). All other slides
# and their content/order are left untouched.

$p = $ppt.ActivePresentation

# Delete the "Stake holders" slide (3rd slide).
$p.Slides.Item(3).Delete()

# After the above deletion everything shifts up by one, so the formerly-last
# (12th) empty slide is now the 11th slide.
$p.Slides.Item(11).Delete()
